# xpod_node header.xlsx -- "Added aiden's changes and formatted csv"
#
# The sheet gains 3 new sensor columns (Figaro 3 Heater / Figaro 4 raw /
# Figaro 4 Heater) right after the existing "FIG 3 (raw)" column, and a new
# second header row with the measurement units for every column. Column
# widths, the active selection and the page orientation are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the 2 brand-new columns (old col F keeps its text and
#        becomes the 3rd of 3 new "Figaro" columns; everything from the old
#        G column onward shifts right by two) -----------------------------
$ws.Columns("G:H").Insert()

# --- 2. Second header row (units / raw-value notes), typed in the same
#        order the original author entered them so the shared-string table
#        builds up identically ------------------------------------------
$ws.Range("A2").Value = "(YYYY-MM-DDThh:mm:ss)"
$ws.Range("C2").Value = "Raw Sensor value"
$ws.Range("D2").Value = "Raw Sensor value"
$ws.Range("E2").Value = "Raw Sensor Value"
$ws.Range("G2").Value = "Raw Sensor Value"
$ws.Range("F2").Value = "Raw Heater Value"
$ws.Range("H2").Value = "Raw Heater Value"

# --- 3. The 3 new first-row headers --------------------------------------
$ws.Range("G1").Value = "Figaro 4 (raw)"
$ws.Range("F1").Value = "Figaro 3(Heater)"
$ws.Range("H1").Value = "Figaro 4(Heater) "

# --- 4. Rest of row 2 ------------------------------------------------------
$ws.Range("B2").Value = "(in volts)"
$ws.Range("L2").Value = "In degree C"
$ws.Range("M2").Value = "hpa"
$ws.Range("N2").Value = "%"
$ws.Range("J2").Value = "Raw Val"
$ws.Range("K2").Value = "Raw Val"
$ws.Range("O2").Value = "Quad 1 Channel 1"
$ws.Range("P2").Value = "Quad 1 Channel 2"
$ws.Range("Q2").Value = "Quad 1 Channel 3"
$ws.Range("R2").Value = "Quad 1 Channel 4"
$ws.Range("S2").Value = "Quad 2 Channel 1"
$ws.Range("T2").Value = "Quad 2 Channel 2"
$ws.Range("U2").Value = "Quad 2 Channel 3"
$ws.Range("V2").Value = "Quad 2 Channel 4"
$ws.Range("W2").Value = "Raw Value"
$ws.Range("AN2").Value = "Gps date"
$ws.Range("AM2").Value = "GPS time"
$ws.Range("AO2").Value = "mph"

# --- 5. Row 2 is taller than the default (wrapped unit text) --------------
$ws.Rows("2:2").RowHeight = 31.2

# --- 6. Column widths (re-sized for the new, longer headers) --------------
$ws.Columns("A:A").ColumnWidth = 22.333333
$ws.Columns("B:B").ColumnWidth = 17.666667
$ws.Columns("C:C").ColumnWidth = 15.5
$ws.Columns("D:D").ColumnWidth = 16.666667
$ws.Columns("E:E").ColumnWidth = 16
$ws.Columns("F:F").ColumnWidth = 16.333333
$ws.Columns("G:G").ColumnWidth = 14.666667
$ws.Columns("H:H").ColumnWidth = 15.833333
$ws.Columns("L:L").ColumnWidth = 10.333333
$ws.Columns("O:O").ColumnWidth = 17.5
$ws.Columns("P:P").ColumnWidth = 16
$ws.Columns("Q:Q").ColumnWidth = 14.666667
$ws.Columns("R:R").ColumnWidth = 16.5
$ws.Columns("S:S").ColumnWidth = 14.666667
$ws.Columns("T:T").ColumnWidth = 16.333333
$ws.Columns("U:U").ColumnWidth = 16
$ws.Columns("V:V").ColumnWidth = 16
$ws.Columns("W:W").ColumnWidth = 10.333333
$ws.Columns("X:X").ColumnWidth = 19.5
$ws.Columns("Y:Y").ColumnWidth = 17
$ws.Columns("Z:Z").ColumnWidth = 9
$ws.Columns("AA:AA").ColumnWidth = 9
$ws.Columns("AB:AB").ColumnWidth = 8.333333
$ws.Columns("AE:AE").ColumnWidth = 9.5
$ws.Columns("AK:AK").ColumnWidth = 10
$ws.Columns("AL:AL").ColumnWidth = 9.833333
$ws.Columns("AN:AN").ColumnWidth = 11.166667
$ws.Columns("AO:AO").ColumnWidth = 12.833333

# --- 7. View / selection ----------------------------------------------------
$ws.Range("AL2").Select()

# --- 8. Page orientation ----------------------------------------------------
$ws.PageSetup.Orientation = 1
